$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values (columns D,L,M,N,O,P,Q,S,T) before overwriting, since the
# edit is a row permutation and targets/sources overlap.
$orig = @{}
$orig[2] = @{
    'D' = 44873
    'L' = 'Primera'
    'M' = 300
    'N' = 22000
    'O' = 22500
    'P' = 22250
    'Q' = '$/bandeja 8 kilos'
    'S' = 2781
    'T' = 8
}
$orig[3] = @{
    'D' = 44495
    'L' = 'Segunda'
    'M' = 270
    'N' = 19000
    'O' = 20000
    'P' = 19556
    'Q' = '$/bandeja 8 kilos'
    'S' = 2444
    'T' = 8
}
$orig[4] = @{
    'D' = 44523
    'L' = 'Primera'
    'M' = 400
    'N' = 21000
    'O' = 22000
    'P' = 21500
    'Q' = '$/bandeja 8 kilos'
    'S' = 2688
    'T' = 8
}
$orig[5] = @{
    'D' = 44523
    'L' = 'Segunda'
    'M' = 100
    'N' = 18000
    'O' = 18000
    'P' = 18000
    'Q' = '$/bandeja 8 kilos'
    'S' = 2250
    'T' = 8
}
$orig[6] = @{
    'D' = 44519
    'L' = 'Primera'
    'M' = 400
    'N' = 21000
    'O' = 22000
    'P' = 21500
    'Q' = '$/bandeja 8 kilos'
    'S' = 2688
    'T' = 8
}
$orig[7] = @{
    'D' = 44519
    'L' = 'Segunda'
    'M' = 200
    'N' = 18000
    'O' = 18000
    'P' = 18000
    'Q' = '$/bandeja 8 kilos'
    'S' = 2250
    'T' = 8
}
$orig[8] = @{
    'D' = 44505
    'L' = 'Segunda'
    'M' = 300
    'N' = 19000
    'O' = 20000
    'P' = 19500
    'Q' = '$/bandeja 8 kilos'
    'S' = 2438
    'T' = 8
}
$orig[9] = @{
    'D' = 44876
    'L' = 'Primera'
    'M' = 300
    'N' = 22000
    'O' = 22500
    'P' = 22250
    'Q' = '$/bandeja 8 kilos'
    'S' = 2781
    'T' = 8
}
$orig[10] = @{
    'D' = 44159
    'L' = 'Primera'
    'M' = 300
    'N' = 2000
    'O' = 2100
    'P' = 2050
    'Q' = '$/kilo (en caja de 14 kilos)'
    'S' = 2050
    'T' = 1
}
$orig[11] = @{
    'D' = 44530
    'L' = 'Primera'
    'M' = 200
    'N' = 19000
    'O' = 20000
    'P' = 19500
    'Q' = '$/bandeja 8 kilos'
    'S' = 2438
    'T' = 8
}
$orig[12] = @{
    'D' = 44530
    'L' = 'Segunda'
    'M' = 100
    'N' = 16000
    'O' = 16000
    'P' = 16000
    'Q' = '$/bandeja 8 kilos'
    'S' = 2000
    'T' = 8
}
$orig[13] = @{
    'D' = 44526
    'L' = 'Primera'
    'M' = 300
    'N' = 21000
    'O' = 21000
    'P' = 21000
    'Q' = '$/bandeja 8 kilos'
    'S' = 2625
    'T' = 8
}
$orig[14] = @{
    'D' = 44488
    'L' = 'Segunda'
    'M' = 160
    'N' = 17000
    'O' = 18000
    'P' = 17500
    'Q' = '$/bandeja 8 kilos'
    'S' = 2188
    'T' = 8
}
$orig[15] = @{
    'D' = 44880
    'L' = 'Primera'
    'M' = 300
    'N' = 22000
    'O' = 22500
    'P' = 22250
    'Q' = '$/bandeja 8 kilos'
    'S' = 2781
    'T' = 8
}
$orig[16] = @{
    'D' = 44512
    'L' = 'Segunda'
    'M' = 300
    'N' = 19000
    'O' = 20000
    'P' = 19500
    'Q' = '$/bandeja 8 kilos'
    'S' = 2438
    'T' = 8
}
$orig[17] = @{
    'D' = 44509
    'L' = 'Segunda'
    'M' = 200
    'N' = 19000
    'O' = 20000
    'P' = 19500
    'Q' = '$/bandeja 8 kilos'
    'S' = 2438
    'T' = 8
}
$orig[18] = @{
    'D' = 44894
    'L' = 'Primera'
    'M' = 200
    'N' = 22000
    'O' = 22500
    'P' = 22250
    'Q' = '$/bandeja 8 kilos'
    'S' = 2781
    'T' = 8
}
$orig[19] = @{
    'D' = 44498
    'L' = 'Segunda'
    'M' = 300
    'N' = 19000
    'O' = 20000
    'P' = 19500
    'Q' = '$/bandeja 8 kilos'
    'S' = 2438
    'T' = 8
}
$orig[20] = @{
    'D' = 44895
    'L' = 'Primera'
    'M' = 200
    'N' = 22000
    'O' = 22500
    'P' = 22250
    'Q' = '$/bandeja 8 kilos'
    'S' = 2781
    'T' = 8
}
$orig[21] = @{
    'D' = 44491
    'L' = 'Segunda'
    'M' = 200
    'N' = 18000
    'O' = 19000
    'P' = 18500
    'Q' = '$/bandeja 8 kilos'
    'S' = 2312
    'T' = 8
}
$orig[22] = @{
    'D' = 44516
    'L' = 'Segunda'
    'M' = 200
    'N' = 18000
    'O' = 19000
    'P' = 18500
    'Q' = '$/bandeja 8 kilos'
    'S' = 2312
    'T' = 8
}
$orig[23] = @{
    'D' = 44162
    'L' = 'Primera'
    'M' = 200
    'N' = 2000
    'O' = 2100
    'P' = 2050
    'Q' = '$/kilo (en caja de 14 kilos)'
    'S' = 2050
    'T' = 1
}
$orig[24] = @{
    'D' = 44533
    'L' = 'Primera'
    'M' = 300
    'N' = 18000
    'O' = 19000
    'P' = 18500
    'Q' = '$/bandeja 8 kilos'
    'S' = 2312
    'T' = 8
}
$orig[25] = @{
    'D' = 44533
    'L' = 'Segunda'
    'M' = 100
    'N' = 16000
    'O' = 16000
    'P' = 16000
    'Q' = '$/bandeja 8 kilos'
    'S' = 2000
    'T' = 8
}
$orig[26] = @{
    'D' = 44890
    'L' = 'Primera'
    'M' = 200
    'N' = 22000
    'O' = 22500
    'P' = 22250
    'Q' = '$/bandeja 8 kilos'
    'S' = 2781
    'T' = 8
}

# Apply permutation: target row <- source row (original values)
$src = $orig[18]
$ws.Range("D2").Value = $src['D']
$ws.Range("L2").Value = $src['L']
$ws.Range("M2").Value = $src['M']
$ws.Range("N2").Value = $src['N']
$ws.Range("O2").Value = $src['O']
$ws.Range("P2").Value = $src['P']
$ws.Range("Q2").Value = $src['Q']
$ws.Range("S2").Value = $src['S']
$ws.Range("T2").Value = $src['T']

$src = $orig[8]
$ws.Range("D3").Value = $src['D']
$ws.Range("L3").Value = $src['L']
$ws.Range("M3").Value = $src['M']
$ws.Range("N3").Value = $src['N']
$ws.Range("O3").Value = $src['O']
$ws.Range("P3").Value = $src['P']
$ws.Range("Q3").Value = $src['Q']
$ws.Range("S3").Value = $src['S']
$ws.Range("T3").Value = $src['T']

$src = $orig[16]
$ws.Range("D4").Value = $src['D']
$ws.Range("L4").Value = $src['L']
$ws.Range("M4").Value = $src['M']
$ws.Range("N4").Value = $src['N']
$ws.Range("O4").Value = $src['O']
$ws.Range("P4").Value = $src['P']
$ws.Range("Q4").Value = $src['Q']
$ws.Range("S4").Value = $src['S']
$ws.Range("T4").Value = $src['T']

$src = $orig[26]
$ws.Range("D5").Value = $src['D']
$ws.Range("L5").Value = $src['L']
$ws.Range("M5").Value = $src['M']
$ws.Range("N5").Value = $src['N']
$ws.Range("O5").Value = $src['O']
$ws.Range("P5").Value = $src['P']
$ws.Range("Q5").Value = $src['Q']
$ws.Range("S5").Value = $src['S']
$ws.Range("T5").Value = $src['T']

$src = $orig[23]
$ws.Range("D6").Value = $src['D']
$ws.Range("L6").Value = $src['L']
$ws.Range("M6").Value = $src['M']
$ws.Range("N6").Value = $src['N']
$ws.Range("O6").Value = $src['O']
$ws.Range("P6").Value = $src['P']
$ws.Range("Q6").Value = $src['Q']
$ws.Range("S6").Value = $src['S']
$ws.Range("T6").Value = $src['T']

$src = $orig[24]
$ws.Range("D7").Value = $src['D']
$ws.Range("L7").Value = $src['L']
$ws.Range("M7").Value = $src['M']
$ws.Range("N7").Value = $src['N']
$ws.Range("O7").Value = $src['O']
$ws.Range("P7").Value = $src['P']
$ws.Range("Q7").Value = $src['Q']
$ws.Range("S7").Value = $src['S']
$ws.Range("T7").Value = $src['T']

$src = $orig[25]
$ws.Range("D8").Value = $src['D']
$ws.Range("L8").Value = $src['L']
$ws.Range("M8").Value = $src['M']
$ws.Range("N8").Value = $src['N']
$ws.Range("O8").Value = $src['O']
$ws.Range("P8").Value = $src['P']
$ws.Range("Q8").Value = $src['Q']
$ws.Range("S8").Value = $src['S']
$ws.Range("T8").Value = $src['T']

$src = $orig[11]
$ws.Range("D9").Value = $src['D']
$ws.Range("L9").Value = $src['L']
$ws.Range("M9").Value = $src['M']
$ws.Range("N9").Value = $src['N']
$ws.Range("O9").Value = $src['O']
$ws.Range("P9").Value = $src['P']
$ws.Range("Q9").Value = $src['Q']
$ws.Range("S9").Value = $src['S']
$ws.Range("T9").Value = $src['T']

$src = $orig[12]
$ws.Range("D10").Value = $src['D']
$ws.Range("L10").Value = $src['L']
$ws.Range("M10").Value = $src['M']
$ws.Range("N10").Value = $src['N']
$ws.Range("O10").Value = $src['O']
$ws.Range("P10").Value = $src['P']
$ws.Range("Q10").Value = $src['Q']
$ws.Range("S10").Value = $src['S']
$ws.Range("T10").Value = $src['T']

$src = $orig[13]
$ws.Range("D11").Value = $src['D']
$ws.Range("L11").Value = $src['L']
$ws.Range("M11").Value = $src['M']
$ws.Range("N11").Value = $src['N']
$ws.Range("O11").Value = $src['O']
$ws.Range("P11").Value = $src['P']
$ws.Range("Q11").Value = $src['Q']
$ws.Range("S11").Value = $src['S']
$ws.Range("T11").Value = $src['T']

$src = $orig[15]
$ws.Range("D12").Value = $src['D']
$ws.Range("L12").Value = $src['L']
$ws.Range("M12").Value = $src['M']
$ws.Range("N12").Value = $src['N']
$ws.Range("O12").Value = $src['O']
$ws.Range("P12").Value = $src['P']
$ws.Range("Q12").Value = $src['Q']
$ws.Range("S12").Value = $src['S']
$ws.Range("T12").Value = $src['T']

$src = $orig[9]
$ws.Range("D13").Value = $src['D']
$ws.Range("L13").Value = $src['L']
$ws.Range("M13").Value = $src['M']
$ws.Range("N13").Value = $src['N']
$ws.Range("O13").Value = $src['O']
$ws.Range("P13").Value = $src['P']
$ws.Range("Q13").Value = $src['Q']
$ws.Range("S13").Value = $src['S']
$ws.Range("T13").Value = $src['T']

$src = $orig[17]
$ws.Range("D14").Value = $src['D']
$ws.Range("L14").Value = $src['L']
$ws.Range("M14").Value = $src['M']
$ws.Range("N14").Value = $src['N']
$ws.Range("O14").Value = $src['O']
$ws.Range("P14").Value = $src['P']
$ws.Range("Q14").Value = $src['Q']
$ws.Range("S14").Value = $src['S']
$ws.Range("T14").Value = $src['T']

$src = $orig[4]
$ws.Range("D15").Value = $src['D']
$ws.Range("L15").Value = $src['L']
$ws.Range("M15").Value = $src['M']
$ws.Range("N15").Value = $src['N']
$ws.Range("O15").Value = $src['O']
$ws.Range("P15").Value = $src['P']
$ws.Range("Q15").Value = $src['Q']
$ws.Range("S15").Value = $src['S']
$ws.Range("T15").Value = $src['T']

$src = $orig[5]
$ws.Range("D16").Value = $src['D']
$ws.Range("L16").Value = $src['L']
$ws.Range("M16").Value = $src['M']
$ws.Range("N16").Value = $src['N']
$ws.Range("O16").Value = $src['O']
$ws.Range("P16").Value = $src['P']
$ws.Range("Q16").Value = $src['Q']
$ws.Range("S16").Value = $src['S']
$ws.Range("T16").Value = $src['T']

$src = $orig[3]
$ws.Range("D17").Value = $src['D']
$ws.Range("L17").Value = $src['L']
$ws.Range("M17").Value = $src['M']
$ws.Range("N17").Value = $src['N']
$ws.Range("O17").Value = $src['O']
$ws.Range("P17").Value = $src['P']
$ws.Range("Q17").Value = $src['Q']
$ws.Range("S17").Value = $src['S']
$ws.Range("T17").Value = $src['T']

$src = $orig[22]
$ws.Range("D18").Value = $src['D']
$ws.Range("L18").Value = $src['L']
$ws.Range("M18").Value = $src['M']
$ws.Range("N18").Value = $src['N']
$ws.Range("O18").Value = $src['O']
$ws.Range("P18").Value = $src['P']
$ws.Range("Q18").Value = $src['Q']
$ws.Range("S18").Value = $src['S']
$ws.Range("T18").Value = $src['T']

$src = $orig[10]
$ws.Range("D19").Value = $src['D']
$ws.Range("L19").Value = $src['L']
$ws.Range("M19").Value = $src['M']
$ws.Range("N19").Value = $src['N']
$ws.Range("O19").Value = $src['O']
$ws.Range("P19").Value = $src['P']
$ws.Range("Q19").Value = $src['Q']
$ws.Range("S19").Value = $src['S']
$ws.Range("T19").Value = $src['T']

$src = $orig[6]
$ws.Range("D20").Value = $src['D']
$ws.Range("L20").Value = $src['L']
$ws.Range("M20").Value = $src['M']
$ws.Range("N20").Value = $src['N']
$ws.Range("O20").Value = $src['O']
$ws.Range("P20").Value = $src['P']
$ws.Range("Q20").Value = $src['Q']
$ws.Range("S20").Value = $src['S']
$ws.Range("T20").Value = $src['T']

$src = $orig[7]
$ws.Range("D21").Value = $src['D']
$ws.Range("L21").Value = $src['L']
$ws.Range("M21").Value = $src['M']
$ws.Range("N21").Value = $src['N']
$ws.Range("O21").Value = $src['O']
$ws.Range("P21").Value = $src['P']
$ws.Range("Q21").Value = $src['Q']
$ws.Range("S21").Value = $src['S']
$ws.Range("T21").Value = $src['T']

$src = $orig[20]
$ws.Range("D22").Value = $src['D']
$ws.Range("L22").Value = $src['L']
$ws.Range("M22").Value = $src['M']
$ws.Range("N22").Value = $src['N']
$ws.Range("O22").Value = $src['O']
$ws.Range("P22").Value = $src['P']
$ws.Range("Q22").Value = $src['Q']
$ws.Range("S22").Value = $src['S']
$ws.Range("T22").Value = $src['T']

$src = $orig[2]
$ws.Range("D23").Value = $src['D']
$ws.Range("L23").Value = $src['L']
$ws.Range("M23").Value = $src['M']
$ws.Range("N23").Value = $src['N']
$ws.Range("O23").Value = $src['O']
$ws.Range("P23").Value = $src['P']
$ws.Range("Q23").Value = $src['Q']
$ws.Range("S23").Value = $src['S']
$ws.Range("T23").Value = $src['T']

$src = $orig[19]
$ws.Range("D24").Value = $src['D']
$ws.Range("L24").Value = $src['L']
$ws.Range("M24").Value = $src['M']
$ws.Range("N24").Value = $src['N']
$ws.Range("O24").Value = $src['O']
$ws.Range("P24").Value = $src['P']
$ws.Range("Q24").Value = $src['Q']
$ws.Range("S24").Value = $src['S']
$ws.Range("T24").Value = $src['T']

$src = $orig[21]
$ws.Range("D25").Value = $src['D']
$ws.Range("L25").Value = $src['L']
$ws.Range("M25").Value = $src['M']
$ws.Range("N25").Value = $src['N']
$ws.Range("O25").Value = $src['O']
$ws.Range("P25").Value = $src['P']
$ws.Range("Q25").Value = $src['Q']
$ws.Range("S25").Value = $src['S']
$ws.Range("T25").Value = $src['T']

$src = $orig[14]
$ws.Range("D26").Value = $src['D']
$ws.Range("L26").Value = $src['L']
$ws.Range("M26").Value = $src['M']
$ws.Range("N26").Value = $src['N']
$ws.Range("O26").Value = $src['O']
$ws.Range("P26").Value = $src['P']
$ws.Range("Q26").Value = $src['Q']
$ws.Range("S26").Value = $src['S']
$ws.Range("T26").Value = $src['T']

Write-Host "Applied permutation edit to rows 2-26"